$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.449.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = "'2.932.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'594.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.82%  '
$ws.Range("D6").Value = "'145.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = "'0.504"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("D9").Value = "'6.94"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.86%  '
$ws.Range("E10").Value = '  -1.25%  '
$ws.Range("D11").Value = "'0.440"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("D12").Value = "'0.0000225"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.71%  '
$ws.Range("D13").Value = "'33.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("E14").Value = '  +0.71%  '
$ws.Range("D15").Value = "'3.421.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.54%  '
$ws.Range("D16").Value = "'61.338.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").Value = "'6.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = "'2.934.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").Value = "'432.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D20").Value = "'13.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.22%  '
$ws.Range("E21").Value = '  -0.59%  '
$ws.Range("D22").Value = "'7.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.74%  '
$ws.Range("D23").Value = "'81.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.11%  '
$ws.Range("D24").Value = "'10.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.19%  '
$ws.Range("D25").Value = "'2.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.22%  '
$ws.Range("D26").Value = "'11.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.86%  '
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").Value = "'2.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.34%  '
$ws.Range("E29").Value = '  -0.18%  '
$ws.Range("D30").Value = "'6.95"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("E31").Value = '  +1.61%  '
$ws.Range("D32").Value = "'26.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("E33").Value = '  +0.11%  '
$ws.Range("D34").Value = "'0.0₃0880"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.66%  '
$ws.Range("D35").Value = "'1.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").Value = "'5.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("D37").Value = "'3.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.37%  '
$ws.Range("D38").Value = "'0.124"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("E40").Value = '  +0.26%  '
$ws.Range("D41").Value = "'42.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.82%  '
$ws.Range("E42").Value = '  -1.80%  '
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("D44").Value = "'2.710.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").Value = "'369.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.13%  '
$ws.Range("D46").Value = "'134.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.72%  '
$ws.Range("D48").Value = "'23.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.94%  '
$ws.Range("D49").Value = "'0.105"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.95%  '
$ws.Range("E50").Value = '  -1.56%  '
$ws.Range("E51").Value = '  -0.30%  '
